$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision of the existing last row's timestamp (A15)
$ws.Range("A15").Value = 45874.50230082176

# Append the new data row (row 16) with the latest sensor reading
$ws.Range("A16").Value = 45874.54183811676
$ws.Range("B16").Value = 2025
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 19.78
$ws.Range("E16").Value = 76.77
$ws.Range("F16").Value = 578.3200000000001
$ws.Range("G16").Value = 14.38
$ws.Range("H16").Value = "ESE"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "13:00:14"

# Match the number format of the date cell in column A to the existing cells above it
$ws.Range("A16").NumberFormat = $ws.Range("A15").NumberFormat
